$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 100002390
$ws.Range("J40").Value = 166669310
$ws.Range("L40").Value = 166669310
$ws.Range("N40").Value = -166669660
# Row 112
$ws.Range("H112").Value = 1890.4482
$ws.Range("J112").Value = 1901
$ws.Range("L112").Value = 5703
$ws.Range("N112").Value = -7919
# Row 124
$ws.Range("H124").Value = 99999
$ws.Range("J124").Value = 99999
$ws.Range("L124").Value = 99999
$ws.Range("N124").Value = -109819
# Row 125
$ws.Range("H125").Value = 15249.75
$ws.Range("I125").Value = 14999.667
$ws.Range("K125").Value = 134997.003
$ws.Range("M125").Value = -132537.003
# Row 132
$ws.Range("H132").Value = 2733.4243
$ws.Range("I132").Value = 2309.7666
$ws.Range("J132").Value = 6970
$ws.Range("K132").Value = 6929.2998
$ws.Range("L132").Value = 20910
$ws.Range("M132").Value = -4399.2998
$ws.Range("N132").Value = -25970

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3634.75
$ws.Range("I45").Value = 2234.8572
$ws.Range("K45").Value = 2234.8572
$ws.Range("M45").Value = -1857.8572
# Row 61
$ws.Range("I61").Value = 70006980
$ws.Range("K61").Value = 70006980
$ws.Range("M61").Value = -70006768
# Row 74
$ws.Range("H74").Value = 2379.111
$ws.Range("I74").Value = 2368.5652
$ws.Range("J74").Value = 2439.75
$ws.Range("K74").Value = 2368.5652
$ws.Range("L74").Value = 2439.75
$ws.Range("M74").Value = -1494.5652
$ws.Range("N74").Value = -4187.75
# Row 77
$ws.Range("H77").Value = 2379.111
$ws.Range("I77").Value = 2368.5652
$ws.Range("J77").Value = 2439.75
$ws.Range("K77").Value = 11842.826
$ws.Range("L77").Value = 12198.75
$ws.Range("M77").Value = -7474.826000000001
$ws.Range("N77").Value = -20934.75
# Row 122
$ws.Range("H122").Value = 3950.5715
$ws.Range("I122").Value = 3817.875
$ws.Range("K122").Value = 11453.625
$ws.Range("M122").Value = -9003.625
# Row 132
$ws.Range("H132").Value = 2003929.9
$ws.Range("I132").Value = 3845.2559
$ws.Range("J132").Value = 14290164
$ws.Range("K132").Value = 11535.7677
$ws.Range("L132").Value = 42870492
$ws.Range("M132").Value = -9005.7677
$ws.Range("N132").Value = -42875552
# Row 136
$ws.Range("I136").Value = 70006980
$ws.Range("K136").Value = 210020940
$ws.Range("M136").Value = -210018390

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 743.6316
$ws.Range("I20").Value = 702.2
$ws.Range("K20").Value = 702.2
$ws.Range("M20").Value = -455.2
# Row 86
$ws.Range("H86").Value = 2990.6897
$ws.Range("I86").Value = 2244.4092
$ws.Range("K86").Value = 2244.4092
$ws.Range("M86").Value = -1121.4092
# Row 89
$ws.Range("H89").Value = 2990.6897
$ws.Range("I89").Value = 2244.4092
$ws.Range("K89").Value = 11222.046
$ws.Range("M89").Value = -5606.046
# Row 95
$ws.Range("H95").Value = 50901
$ws.Range("J95").Value = 50901
$ws.Range("L95").Value = 50901
$ws.Range("N95").Value = -56393
# Row 107
$ws.Range("H107").Value = 5883
$ws.Range("I107").Value = 6009.143
$ws.Range("K107").Value = 6009.143
$ws.Range("M107").Value = -4089.143
# Row 134
$ws.Range("H134").Value = 3848222.5
$ws.Range("I134").Value = 2120.25
$ws.Range("J134").Value = 50001450
$ws.Range("K134").Value = 6360.75
$ws.Range("L134").Value = 150004350
$ws.Range("M134").Value = -3825.75
$ws.Range("N134").Value = -150009420

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 27781448
$ws.Range("J31").Value = 4600.4614
$ws.Range("L31").Value = 4600.4614
$ws.Range("N31").Value = -5190.4614
# Row 34
$ws.Range("H34").Value = 27781448
$ws.Range("J34").Value = 4600.4614
$ws.Range("L34").Value = 4600.4614
$ws.Range("N34").Value = -5004.4614
# Row 58
$ws.Range("H58").Value = 2674.0386
$ws.Range("I58").Value = 2499.7778
$ws.Range("K58").Value = 2499.7778
$ws.Range("M58").Value = -2296.7778
# Row 59
$ws.Range("H59").Value = 124999.5
$ws.Range("J59").Value = 124999.5
$ws.Range("L59").Value = 124999.5
$ws.Range("N59").Value = -127289.5
# Row 68
$ws.Range("H68").Value = 99999
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
# Row 71
$ws.Range("H71").Value = 99999
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
# Row 103
$ws.Range("H103").Value = 32182.572
$ws.Range("I103").Value = 11056
$ws.Range("K103").Value = 11056
$ws.Range("M103").Value = -9884
# Row 136
$ws.Range("H136").Value = 2674.0386
$ws.Range("I136").Value = 2499.7778
$ws.Range("K136").Value = 7499.3334
$ws.Range("M136").Value = -4949.3334

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 465.25
$ws.Range("I5").Value = 354.125
$ws.Range("K5").Value = 1062.375
$ws.Range("M5").Value = -950.375
# Row 62
$ws.Range("H62").Value = 13672.167
$ws.Range("I62").Value = 7911.5
$ws.Range("J62").Value = 14824.3
$ws.Range("K62").Value = 23734.5
$ws.Range("L62").Value = 44472.89999999999
$ws.Range("M62").Value = -23048.5
$ws.Range("N62").Value = -45844.89999999999
# Row 63
$ws.Range("H63").Value = 20988.066
$ws.Range("I63").Value = 13165.667
$ws.Range("K63").Value = 39497.001
$ws.Range("M63").Value = -38748.001
# Row 65
$ws.Range("H65").Value = 13672.167
$ws.Range("I65").Value = 7911.5
$ws.Range("J65").Value = 14824.3
$ws.Range("K65").Value = 71203.5
$ws.Range("L65").Value = 133418.7
$ws.Range("M65").Value = -67771.5
$ws.Range("N65").Value = -140282.7
# Row 66
$ws.Range("H66").Value = 20988.066
$ws.Range("I66").Value = 13165.667
$ws.Range("K66").Value = 118491.003
$ws.Range("M66").Value = -114747.003
# Row 69
$ws.Range("H69").Value = 16317
$ws.Range("J69").Value = 16317
$ws.Range("L69").Value = 48951
$ws.Range("N69").Value = -50573
# Row 72
$ws.Range("H72").Value = 16317
$ws.Range("J72").Value = 16317
$ws.Range("L72").Value = 146853
$ws.Range("N72").Value = -154965
# Row 135
$ws.Range("H135").Value = 465.25
$ws.Range("I135").Value = 354.125
$ws.Range("K135").Value = 3187.125
$ws.Range("M135").Value = -652.125

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
# Row 98
$ws.Range("H98").Value = 25820.5
$ws.Range("J98").Value = 25820.5
$ws.Range("L98").Value = 25820.5
$ws.Range("N98").Value = -31810.5
# Row 113
$ws.Range("H113").Value = 1687231.2
$ws.Range("I113").Value = 3785.125
$ws.Range("K113").Value = 3785.125
$ws.Range("M113").Value = -1615.125
# Row 132
$ws.Range("H132").Value = 6821019
$ws.Range("I132").Value = 2977.76
$ws.Range("K132").Value = 8933.280000000001
$ws.Range("M132").Value = -6403.280000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1580.7391
$ws.Range("J16").Value = 2251
$ws.Range("L16").Value = 2251
$ws.Range("N16").Value = -2591
# Row 55
$ws.Range("H55").Value = 994
$ws.Range("J55").Value = 1314.1428
$ws.Range("L55").Value = 1314.1428
$ws.Range("N55").Value = -1660.1428
# Row 68
$ws.Range("H68").Value = 10417841
$ws.Range("I68").Value = 10417841
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 10417841
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -10417092
$ws.Range("N68").ClearContents()
# Row 71
$ws.Range("H71").Value = 10417841
$ws.Range("I71").Value = 10417841
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 52089205
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -52085461
$ws.Range("N71").ClearContents()
# Row 100
$ws.Range("H100").Value = 10430116
$ws.Range("I100").Value = 2599.7778
$ws.Range("K100").Value = 2599.7778
$ws.Range("M100").Value = -2058.7778
# Row 101
$ws.Range("H101").Value = 27286.637
$ws.Range("J101").Value = 27286.637
$ws.Range("L101").Value = 27286.637
$ws.Range("N101").Value = -33776.637
# Row 136
$ws.Range("H136").Value = 1915.1143
$ws.Range("I136").Value = 1020.88464
$ws.Range("K136").Value = 3062.65392
$ws.Range("M136").Value = -512.6539199999997

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Range("H18").Value = 19967.6
$ws.Range("J18").Value = 19970
$ws.Range("L18").Value = 19970
$ws.Range("N18").Value = -20316
# Row 113
$ws.Range("H113").Value = 711.0833
$ws.Range("I113").Value = 703.7778
$ws.Range("J113").Value = 733
$ws.Range("K113").Value = 2111.3334
$ws.Range("L113").Value = 2199
$ws.Range("M113").Value = 58.66660000000002
$ws.Range("N113").Value = -6539
# Row 136
$ws.Range("H136").Value = 238957.45
$ws.Range("I136").Value = 844.8333
$ws.Range("J136").Value = 1667633.1
$ws.Range("K136").Value = 2534.4999
$ws.Range("L136").Value = 5002899.300000001
$ws.Range("M136").Value = 15.5001000000002
$ws.Range("N136").Value = -5007999.300000001
